$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: row 4 corresponds to the
# 375c5b4f-20f5-49eb-b73c-b2cf0be8cb0e...zh-cn.xlf entry.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-24 06:49:46"
$wsZhCn.Range("H4").Value = "2016-03-24 06:50:15"

# "de-de" worksheet: row 4 corresponds to the
# 375c5b4f-20f5-49eb-b73c-b2cf0be8cb0e...de-de.xlf entry.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-24 06:49:50"
$wsDeDe.Range("H4").Value = "2016-03-24 06:50:22"
